$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row: "_old" columns become "_FV2210", "_new" columns
#    become "_FV2304". Column K ("diff") is left untouched.
# ---------------------------------------------------------------------------
$headersFV2210 = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$headersFV2304 = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headersFV2210.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2210[$i]
}

# Column 11 is "diff" and stays the same.

for ($i = 0; $i -lt $headersFV2304.Count; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $headersFV2304[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the data range A1:U88 into an Excel Table ("Table1").
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U88")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# ---------------------------------------------------------------------------
# 3) Freeze the header row (row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
